# The document contains a nested table (the "Stt / Ten tai lieu / Ma tai
# lieu / Noi dung yeu cau / Ly do" table) embedded inside a merged cell of
# the outer form table. This host's Word table-navigation (Tables,
# Range.Tables, Cell.Tables, Columns, ...) only ever resolves to the
# outermost table object here, so the nested table's grid/cell widths
# can't be reached through Table/Column/Cell APIs. The only path that
# exposes the nested table's markup is Range.WordOpenXML / Range.InsertXML
# (flattened OPC package XML), so use that to adjust the two affected
# column widths directly.
$d = $word.ActiveDocument

$xml = $d.Content.WordOpenXML

# Pull out only the /word/document.xml package part and edit it in
# isolation, then feed just that part back through InsertXML. (Round-
# tripping the *entire* flattened package re-mints relationship ids for
# other parts such as the header image, which would corrupt them; editing
# only the document part avoids touching anything outside it.)
$partStart = '<pkg:part pkg:name="/word/document.xml"'
$startIdx = $xml.IndexOf($partStart)
$partEndTag = '</pkg:part>'
$endIdx = $xml.IndexOf($partEndTag, $startIdx) + $partEndTag.Length
$docPart = $xml.Substring($startIdx, $endIdx - $startIdx)

# Column 1 of the nested table: 1843 dxa -> 709 dxa (tblGrid + every tcW)
$docPart = $docPart.Replace('w:w="1843"', 'w:w="709"')
# Column 2 of the nested table: 1559 dxa -> 2693 dxa (tblGrid + every tcW)
$docPart = $docPart.Replace('w:w="1559"', 'w:w="2693"')

$newPkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + $docPart + '</pkg:package>'

$d.Content.InsertXML($newPkg)
